# Update the data table on Sheet1 and move the active selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (spreadsheet row 3) value updates: E3 5->6, G3 2->3, H3 14->13
$ws.Range("E3").Value = 6
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 13

# Move the selection/active cell from I1 to E3, as captured in the saved view state.
$ws.Range("E3").Select()
